$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Authentication switched from raw request body (reqres.in) to Bearer-token ---
# --- against the bookstore API; the "Body" column (D) is no longer needed.   ---

# 1. Re-point the EndPoint column (B) at the new host. Do this before the
#    column deletion so the shared-string table only ever grows/shrinks once.
# First, mutate D3's own format (it is the sole owner of its style) so the
# wrap-text alignment becomes a plain left/vertical-center alignment -- this
# edits the existing cell style in place instead of minting a new one.
$ws.Range("D3").HorizontalAlignment = -4131   # xlLeft
$ws.Range("D3").WrapText = $false

# Copy that now-updated format onto B2/B3 (both currently unstyled / border-only)
# so they end up sharing the very same style entry.
$ws.Range("D3").Copy()
$ws.Range("B2").PasteSpecial(-4122)           # xlPasteFormats
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. New endpoint URL for both the Get and Post rows.
$ws.Range("B2").Value2 = "https://bookstore.toolsqa.com"
$ws.Range("B3").Value2 = "https://bookstore.toolsqa.com"

# 3. Drop the old "Body" column (D) entirely -- no more raw JSON payload.
$ws.Columns.Item(4).Delete()

# Row 3 no longer needs the tall 60pt height that fit the wrapped JSON body.
$ws.Rows.Item(3).AutoFit()

# 4. Column B needs to widen slightly to fit the new (longer) URL text.
$ws.Columns.Item(2).ColumnWidth = 28

# 5. Selection moves to C2.
[void]$ws.Range("C2").Select()
